# Commit: "fixed export and fixing maps"
#
# - Rename the single worksheet from the generic "1" to the municipality
#   name "საგარეჯო".
# - Drop the "(მოსახლეობის აღწერის შედეგებით)" caption row (old row 2) —
#   it is no longer part of the simplified export.
# - Drop the obsolete 1989 and 2002 census columns, keeping only the 2014
#   area figure.
# - Leave the active selection on A2 (now the blank spacer row under the
#   title), matching where the author's cursor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "საგარეჯო"

# Remove the caption row "(მოსახლეობის აღწერის შედეგებით)" (old row 2).
$ws.Range("A2").EntireRow.Delete()

# Remove the 1989 (B) and 2002 (C) columns, leaving only 2014 (old D, now B).
$ws.Range("B1:C1").EntireColumn.Delete()

$ws.Range("A2").Select()
